# Reorder the "Recorded By" (column G) comma-separated list in every row:
# move the last item of the list to the front (right-rotate the list by 1).
# Cells with a single value (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($val -and $val.ToString().Contains(",")) {
        $rawParts = $val.ToString().Split(",")
        $parts = @()
        foreach ($p in $rawParts) {
            $parts += $p.Trim()
        }

        $count = $parts.Length
        $last = $parts[$count - 1]

        $newParts = @($last)
        for ($i = 0; $i -lt ($count - 1); $i++) {
            $newParts += $parts[$i]
        }

        $newVal = [string]::Join(", ", $newParts)
        $cell.Value = $newVal
    }
}
